# Updated cryptos list (price/volume refresh + a couple of row swaps) as
# published by the GitHub Actions job.
#
# Note: several "Price" values (column D) look numeric (e.g. "0.9100",
# "1.000") but must stay stored as literal text, exactly like the source
# workbook does (t="inlineStr"). A leading single-quote is Excel's classic
# "force text / keep leading-trailing zeros" prefix, so values below are
# written as '<value> (e.g. '''0.9100' -> the literal text 0.9100).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.478.56'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '''1.880.86'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''0.7198'
$ws.Range("E5").Value = '  +2.00%  '
$ws.Range("D6").Value = '''243.05'
$ws.Range("E6").Value = '  +2.21%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.07884'
$ws.Range("E8").Value = '  -1.82%  '
$ws.Range("D9").Value = '''0.3131'
$ws.Range("E9").Value = '  +3.63%  '
$ws.Range("D10").Value = '''25.26'
$ws.Range("E10").Value = '  +7.63%  '
$ws.Range("D11").Value = '''0.08254'
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = '''1.896.14'
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("D13").Value = '''0.7333'
$ws.Range("E13").Value = '  +4.20%  '
$ws.Range("D14").Value = '''5.304'
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").Value = '''91.23'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '''29.488.18'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").Value = '''5.963'
$ws.Range("E17").Value = '  +2.83%  '
$ws.Range("D18").Value = '''247.80'
$ws.Range("E18").Value = '  +4.33%  '
$ws.Range("D19").Value = '''0.000007913'
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").Value = '''13.35'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '''0.9988'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '''7.992'
$ws.Range("E22").Value = '  +6.94%  '
$ws.Range("D23").Value = '''0.9996'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '''0.1581'
$ws.Range("E24").Value = '  +10.49%  '
$ws.Range("D25").Value = '''164.15'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '''9.078'
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").Value = '''18.36'
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("D28").Value = '''1.363'
$ws.Range("E28").Value = '  -4.06%  '
$ws.Range("D29").Value = '''1.500'
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("D30").Value = '''4.389'
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").Value = '''4.158'
$ws.Range("E31").Value = '  +3.64%  '
$ws.Range("D32").Value = '''0.05304'
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("D33").Value = '''1.939'
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("D35").Value = '''0.7255'
$ws.Range("E35").Value = '  +1.60%  '
$ws.Range("D36").Value = '''2.677'
$ws.Range("E36").Value = '  +0.99%  '
$ws.Range("D37").Value = '''0.01872'
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("D38").Value = '''1.267.27'
$ws.Range("E38").Value = '  +12.20%  '
$ws.Range("D39").Value = '''2.731'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").Value = '''0.9100'
$ws.Range("E40").Value = '  -3.43%  '
$ws.Range("D41").Value = '''74.11'
$ws.Range("E41").Value = '  +5.85%  '
$ws.Range("D42").Value = '''6.113'
$ws.Range("E42").Value = '  +3.01%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''1.000'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''103.78'
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").Value = '''0.5322'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("B46").Value = 'SynthetixNetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D46").Value = '''2.950'
$ws.Range("E46").Value = '  +14.41%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.776'
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '''0.4343'
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").Value = '''9.280'
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("D51").Value = '''7.088'
$ws.Range("E51").Value = '  +2.31%  '
